$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-30 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-01 Monday", 2) | Out-Null
$d.Content.Find.Execute("663÷9=73, 6", $true, $false, $false, $false, $false, $true, 1, $false, "206÷7=29, 3", 2) | Out-Null
$d.Content.Find.Execute("283÷4=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "206÷8=25, 6", 2) | Out-Null
$d.Content.Find.Execute("933÷2=466, 1", $true, $false, $false, $false, $false, $true, 1, $false, "821÷6=136, 5", 2) | Out-Null
$d.Content.Find.Execute("703÷7=100, 3", $true, $false, $false, $false, $false, $true, 1, $false, "340÷7=48, 4", 2) | Out-Null
$d.Content.Find.Execute("964÷8=120, 4", $true, $false, $false, $false, $false, $true, 1, $false, "394÷6=65, 4", 2) | Out-Null
$d.Content.Find.Execute("918÷7=131, 1", $true, $false, $false, $false, $false, $true, 1, $false, "177÷5=35, 2", 2) | Out-Null
$d.Content.Find.Execute("707÷2=353, 1", $true, $false, $false, $false, $false, $true, 1, $false, "570÷9=63, 3", 2) | Out-Null
$d.Content.Find.Execute("261÷9=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "265÷4=66, 1", 2) | Out-Null
$d.Content.Find.Execute("722÷9=80, 2", $true, $false, $false, $false, $false, $true, 1, $false, "663÷2=331, 1", 2) | Out-Null
$d.Content.Find.Execute("935÷5=187, 0", $true, $false, $false, $false, $false, $true, 1, $false, "168÷7=24, 0", 2) | Out-Null
$d.Content.Find.Execute("969÷8=121, 1", $true, $false, $false, $false, $false, $true, 1, $false, "775÷5=155, 0", 2) | Out-Null
$d.Content.Find.Execute("125÷3=41, 2", $true, $false, $false, $false, $false, $true, 1, $false, "706÷6=117, 4", 2) | Out-Null
$d.Content.Find.Execute("556÷6=92, 4", $true, $false, $false, $false, $false, $true, 1, $false, "279÷6=46, 3", 2) | Out-Null
$d.Content.Find.Execute("958÷7=136, 6", $true, $false, $false, $false, $false, $true, 1, $false, "519÷2=259, 1", 2) | Out-Null
$d.Content.Find.Execute("773÷4=193, 1", $true, $false, $false, $false, $false, $true, 1, $false, "984÷5=196, 4", 2) | Out-Null
$d.Content.Find.Execute("423÷8=52, 7", $true, $false, $false, $false, $false, $true, 1, $false, "120÷2=60, 0", 2) | Out-Null
$d.Content.Find.Execute("295÷3=98, 1", $true, $false, $false, $false, $false, $true, 1, $false, "635÷6=105, 5", 2) | Out-Null
$d.Content.Find.Execute("839÷8=104, 7", $true, $false, $false, $false, $false, $true, 1, $false, "782÷2=391, 0", 2) | Out-Null
$d.Content.Find.Execute("113÷3=37, 2", $true, $false, $false, $false, $false, $true, 1, $false, "807÷6=134, 3", 2) | Out-Null
$d.Content.Find.Execute("606÷4=151, 2", $true, $false, $false, $false, $false, $true, 1, $false, "929÷5=185, 4", 2) | Out-Null
$d.Content.Find.Execute("892÷2=446, 0", $true, $false, $false, $false, $false, $true, 1, $false, "120÷6=20, 0", 2) | Out-Null
$d.Content.Find.Execute("319÷8=39, 7", $true, $false, $false, $false, $false, $true, 1, $false, "101÷8=12, 5", 2) | Out-Null
$d.Content.Find.Execute("307÷2=153, 1", $true, $false, $false, $false, $false, $true, 1, $false, "948÷2=474, 0", 2) | Out-Null
$d.Content.Find.Execute("181÷7=25, 6", $true, $false, $false, $false, $false, $true, 1, $false, "933÷9=103, 6", 2) | Out-Null
$d.Content.Find.Execute("184÷6=30, 4", $true, $false, $false, $false, $false, $true, 1, $false, "104÷6=17, 2", 2) | Out-Null
